# Apply the coinranking price/volume refresh described by the commit.
# For the "Price" column (D), several new values are plain decimal numbers
# (e.g. "320.08"). The source data stores these as literal TEXT (t="inlineStr"),
# not numbers, so we force text entry the same way a user would in Excel: a
# leading apostrophe. We then reset the cell's Style back to Normal so the
# quote-prefix flag does not leave a stray style on the cell (matching the
# original, unstyled cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.019.04"
$ws.Range("E2").Value = "  -0.70%  "

$ws.Range("D3").Value = "1.916.58"
$ws.Range("E3").Value = "  +0.16%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'320.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.79%  "

$ws.Range("E6").Value = "  +0.08%  "

$ws.Range("D7").Value = "'0.5034"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.12%  "

$ws.Range("D8").Value = "'0.4026"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.63%  "

$ws.Range("D9").Value = "'0.08253"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.88%  "

$ws.Range("D10").Value = "'1.110"
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").Value = "'42.06"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.91%  "

$ws.Range("D12").Value = "'23.78"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.87%  "

$ws.Range("D13").Value = "1.918.07"
$ws.Range("E13").Value = "  +0.21%  "

$ws.Range("D14").Value = "'6.427"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.02%  "

$ws.Range("D15").Value = "'7.310"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.03%  "

$ws.Range("D16").Value = "'1.001"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.07%  "

$ws.Range("D17").Value = "'92.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.94%  "

$ws.Range("E18").Value = "  -1.55%  "

$ws.Range("D19").Value = "'0.06497"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.94%  "

$ws.Range("D20").Value = "'18.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.47%  "

$ws.Range("D22").Value = "'5.946"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.20%  "

$ws.Range("D23").Value = "30.057.28"
$ws.Range("E23").Value = "  -0.60%  "

$ws.Range("D24").Value = "'11.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.81%  "

$ws.Range("D25").Value = "'2.193"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.48%  "

$ws.Range("D26").Value = "'22.31"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.92%  "

$ws.Range("D27").Value = "2.135.72"
$ws.Range("E27").Value = "  +0.02%  "

$ws.Range("D28").Value = "'161.95"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.20%  "

$ws.Range("D29").Value = "'2.302"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.05%  "

$ws.Range("D30").Value = "'128.94"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.47%  "

$ws.Range("D31").Value = "'1.135"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.83%  "

$ws.Range("D32").Value = "'0.1039"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.43%  "

$ws.Range("D33").Value = "'6.014"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.10%  "

$ws.Range("D34").Value = "'3.822"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.83%  "

$ws.Range("D35").Value = "'0.02445"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.68%  "

$ws.Range("D36").Value = "'5.362"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.13%  "

$ws.Range("D37").Value = "'0.06428"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.30%  "

$ws.Range("D38").Value = "'8.968"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.43%  "

$ws.Range("D39").Value = "'0.2166"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.97%  "

$ws.Range("D40").Value = "'1.207"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.93%  "

$ws.Range("D41").Value = "'0.6433"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.65%  "

$ws.Range("D42").Value = "'11.42"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.04%  "

$ws.Range("D43").Value = "'1.219"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.51%  "

$ws.Range("E44").Value = "  +0.11%  "

$ws.Range("D45").Value = "'13.33"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.35%  "

$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "'2.179"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.08%  "

$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "'0.6005"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.32%  "

$ws.Range("D48").Value = "'3.637"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.74%  "

$ws.Range("D49").Value = "'123.39"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.44%  "

$ws.Range("D50").Value = "'1.216"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.46%  "

$ws.Range("D51").Value = "'78.92"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.74%  "
